# Merge first_name/last_name into a single "name" field, rename/repurpose
# remaining columns, add a clickable mailto hyperlink for the contact email,
# and reapply the sheet's look (column widths + centered alignment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the first_name column (B). Excel shifts last_name, skills,
#        experience, contact one column to the left automatically. ---
$ws.Columns.Item(2).Delete()

# --- 2. Re-purpose the (now shifted) columns: last_name -> name, and
#        update the sample row with the merged full name + new email. ---
$ws.Range("B1").Value = "name"
$ws.Range("B2").Value = "Danny Carey"
$ws.Range("E2").Value = "lateralus@buzzfeed.com"

# --- 3. Reapply column widths (matches the widths the sheet carried for
#        id/name/skills/experience/contact + the now-empty trailing column). ---
$ws.Columns.Item(2).ColumnWidth = 56.85546875
$ws.Columns.Item(3).ColumnWidth = 40
$ws.Columns.Item(4).ColumnWidth = 19.5703125
$ws.Columns.Item(5).ColumnWidth = 37.140625
$ws.Columns.Item(6).ColumnWidth = 42.7109375

# --- 4. Centre-align the header + data rows. ---
$ws.Range("A1:E2").HorizontalAlignment = -4108

# --- 5. Stray formatted (but empty) cell left behind at C8, centred both
#        ways, as in the source file. ---
$ws.Range("C8").HorizontalAlignment = -4108
$ws.Range("C8").VerticalAlignment = -4108

# --- 6. Turn the contact e-mail into a real mailto: hyperlink, centred. ---
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:lateralus@buzzfeed.com")
$ws.Range("E2").HorizontalAlignment = -4108

$ws.Range("E13").Select()
